$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store plain text (e.g. "229.37", "+2.44%") rather
# than numbers. Force the target range to Text format before writing so the
# COM layer doesn't reinterpret numeric-looking strings as numbers, then
# restore the default "Normal" style so no stray formatting is introduced.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '38.658.40'
$ws.Range("E2").Value = '  +2.44%  '
$ws.Range("D3").Value = '2.100.26'
$ws.Range("E3").Value = '  +3.43%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '229.37'
$ws.Range("E5").Value = '  +0.83%  '
$ws.Range("D6").Value = '0.616'
$ws.Range("E6").Value = '  +2.10%  '
$ws.Range("D7").Value = '61.45'
$ws.Range("E7").Value = '  +2.67%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.383'
$ws.Range("E9").Value = '  +1.85%  '
$ws.Range("D10").Value = '0.0844'
$ws.Range("E10").Value = '  +2.21%  '
$ws.Range("D11").Value = '0.105'
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("D12").Value = '2.408.38'
$ws.Range("E12").Value = '  +3.35%  '
$ws.Range("D13").Value = '14.82'
$ws.Range("D14").Value = '22.50'
$ws.Range("E14").Value = '  +7.19%  '
$ws.Range("D15").Value = '0.784'
$ws.Range("E15").Value = '  +2.38%  '
$ws.Range("D16").Value = '5.49'
$ws.Range("E16").Value = '  +5.79%  '
$ws.Range("D17").Value = '2.101.60'
$ws.Range("E17").Value = '  +3.05%  '
$ws.Range("D18").Value = '38.525.40'
$ws.Range("E18").Value = '  +2.22%  '
$ws.Range("D19").Value = '70.95'
$ws.Range("E19").Value = '  +2.11%  '
$ws.Range("D20").Value = '6.04'
$ws.Range("E20").Value = '  +2.46%  '
$ws.Range("D21").Value = '0.0₃0836'
$ws.Range("E21").Value = '  +1.36%  '
$ws.Range("D22").Value = '227.22'
$ws.Range("E22").Value = '  +1.63%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  +3.08%  '
$ws.Range("E25").Value = '  +2.76%  '
$ws.Range("D26").Value = '170.35'
$ws.Range("E26").Value = '  +1.33%  '
$ws.Range("D27").Value = '9.44'
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("E28").Value = '  +1.52%  '
$ws.Range("D29").Value = '19.14'
$ws.Range("E29").Value = '  +1.97%  '
$ws.Range("E30").Value = '  +8.16%  '
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("E32").Value = '  +4.57%  '
$ws.Range("D33").Value = '4.79'
$ws.Range("E33").Value = '  +6.91%  '
$ws.Range("D34").Value = '4.48'
$ws.Range("E34").Value = '  +2.23%  '
$ws.Range("D35").Value = '0.0606'
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("D36").Value = '6.52'
$ws.Range("E36").Value = '  +0.89%  '
$ws.Range("E37").Value = '  +4.72%  '
$ws.Range("D38").Value = '3.57'
$ws.Range("E38").Value = '  +4.15%  '
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("D40").Value = '18.77'
$ws.Range("E40").Value = '  +3.15%  '
$ws.Range("D41").Value = '1.546.41'
$ws.Range("E41").Value = '  +0.48%  '
$ws.Range("D42").Value = '100.07'
$ws.Range("E42").Value = '  +4.69%  '
$ws.Range("E43").Value = '  +2.84%  '
$ws.Range("E44").Value = '  +1.47%  '
$ws.Range("D45").Value = '0.0916'
$ws.Range("E45").Value = '  +1.04%  '
$ws.Range("D46").Value = '4.20'
$ws.Range("E46").Value = '  +3.22%  '
$ws.Range("D47").Value = '7.59'
$ws.Range("E47").Value = '  +6.73%  '
$ws.Range("E48").Value = '  +1.43%  '
$ws.Range("E49").Value = '  +3.71%  '
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("D51").Value = '2.293.85'
$ws.Range("E51").Value = '  +3.35%  '

$ws.Range("D2:E51").Style = "Normal"
